$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new row at position 12 (shifts rows 12-20 down to 13-21) ---
$ws.Rows.Item(12).Insert()

# Copy the formatting (styles/borders) from row 11 (B:D) into the newly
# inserted row 12 (B:D) so it matches its siblings, then set its height.
$ws.Range("B11:D11").Copy()
$ws.Range("B12:D12").PasteSpecial(-4122)
$ws.Rows.Item(12).RowHeight = 19.5

# --- Update text content (renumbering the use-case steps) ---
$ws.Range("D9").Value = "3. Obtém modelos disponiveis"
$ws.Range("D12").Value = "6. Regista modelo escolhido"
$ws.Range("D13").Value = "7. Pergunta qual a forma de personalizar o carro"
$ws.Range("C14").Value = "8. Escolhe configuração ótima"
$ws.Range("D15").Value = "9.<<include>> Escolher configuração ótima"
$ws.Range("B16").Value = "Alternativa 1 [Escolher Pacote] (passo 8)"
$ws.Range("C16").Value = "8.1 Escolhe Pacote"
$ws.Range("D17").Value = "8.2 <<include>> Escolhe Pacote"
$ws.Range("B19").Value = "Alternativa 2 [Escolher Especificações] (passo 8)"
$ws.Range("C19").Value = "8.1 Escolhe Especificações"
$ws.Range("D20").Value = "8.2 <<include>> Escolher Especificações"

# --- Update the view (scroll position / selection) ---
$ws.Application.ActiveWindow.ScrollRow = 6
$ws.Range("C17").Select()

Write-Host "edit complete"
